# The diff updates the "Förändrad" (changed) date in column C for every
# data row (rows 2-216) from serial date 45186 (2023-09-17) to
# serial date 45188 (2023-09-19). No other cells are affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C216").Value = 45188
